$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the hidden _FilterDatabase defined name range (E70 -> E66)
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "sheet1!_FilterDatabase") {
        $n.RefersTo = "=sheet1!`$A`$1:`$E`$66"
    }
}

# ---------------------------------------------------------------------------
# 2. Update rows 2-16 with the new charging-station data
#    (cell styles are left untouched; only the values change)
# ---------------------------------------------------------------------------
$stationNames = @(
    "飞狐四方坪西区充电站",
    "飞狐四方坪南区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪东区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站",
    "飞狐四方坪西区充电站"
)

$terminalNames = @(
    "9176699355900102",
    "9176699368200305",
    "9176699400500202",
    "9176699442100102",
    "9176699400501301",
    "9176699400500404",
    "9176699400500804",
    "9176699400501201",
    "9176699400501101",
    "9176699400501102",
    "9176699400501205",
    "9176699400500104",
    "9176699400500304",
    "9176699400500203",
    "9176699400500701"
)

$lastChargeEnd = @(
    "46025.218819444446",
    "46026.070289351854",
    "46026.569074074076",
    "46026.813402777778",
    "46027.052037037036",
    "46027.097581018519",
    "46027.145127314812",
    "46027.152650462966",
    "46027.177349537036",
    "46027.244513888887",
    "46027.276018518518",
    "46027.282187500001",
    "46027.38045138889",
    "46027.390960648147",
    "46027.394409722219"
)

for ($i = 0; $i -lt $stationNames.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $stationNames[$i]
    $ws.Cells.Item($r, 2).Value = $terminalNames[$i]
    $ws.Cells.Item($r, 3).Value = [double]$lastChargeEnd[$i]
    $ws.Cells.Item($r, 4).Value = 46027.999918981484
}

# ---------------------------------------------------------------------------
# 3. Rows 17-20 no longer hold data -> reset them to the same blank style
#    used by the surrounding empty rows (e.g. row 21) and clear the values
# ---------------------------------------------------------------------------
$ws.Range("A21:E21").Copy()
$ws.Range("A17:E20").PasteSpecial(-4122)
$ws.Range("A17:E20").ClearContents()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Swap the formatting of rows 37 and 41
# ---------------------------------------------------------------------------
$ws.Range("A37:E37").Copy()
$ws.Range("A200").PasteSpecial(-4122)

$ws.Range("A41:E41").Copy()
$ws.Range("A37").PasteSpecial(-4122)

$ws.Range("A200:E200").Copy()
$ws.Range("A41").PasteSpecial(-4122)

$ws.Rows("200:200").Delete()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Rows 52-55, column E: switch style from s=10 to s=5 (same as row 56)
# ---------------------------------------------------------------------------
$ws.Range("E56").Copy()
$ws.Range("E52:E55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Remove rows 61-64 entirely (dimension shrinks from E64 to E60)
# ---------------------------------------------------------------------------
$ws.Rows("61:64").Delete()

# ---------------------------------------------------------------------------
# 7. Update the active selection shown when the sheet is reopened
# ---------------------------------------------------------------------------
$ws.Range("D27").Select()
